$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.569.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "'1.920.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.60%  "
$ws.Range("D5").Value = "'326.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "'0.4818"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'0.4067"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "'0.08245"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").Value = "'1.011"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Value = "'23.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "'1.927.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").Value = "'6.073"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "'7.245"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").Value = "'91.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").Value = "'0.06887"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "'0.00001041"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'17.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "'1.012"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "'29.584.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'5.687"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").Value = "'11.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").Value = "'2.187"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").Value = "'2.150.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'6.554"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.64%  "
$ws.Range("D27").Value = "'155.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").Value = "'19.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "'120.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").Value = "'1.022"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "'0.09638"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").Value = "'5.633"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'1.378"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("D36").Value = "'0.06384"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.45%  "
$ws.Range("D37").Value = "'0.02292"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").Value = "'1.191"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("D39").Value = "'0.5951"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "'10.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").Value = "'7.908"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "'0.1852"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").Value = "'2.468"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "'12.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("D46").Value = "'0.07507"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.92%  "
$ws.Range("D47").Value = "'0.5572"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "'1.945"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("D49").Value = "'118.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("D50").Value = "'2.441"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.80%  "
$ws.Range("D51").Value = "'72.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.63%  "
